$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove column M (old "M" data is discarded; old column N shifts left to
# become the new column M), matching the "remove column from alcohol data"
# commit.
$ws.Range("M1").EntireColumn.Delete()

# Leave the selection where the deleted column used to be.
$ws.Range("M1").Select()
